$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.694027185440063
$ws.Range("B1").Value = 2.003679275512695
$ws.Range("C1").Value = 5.261378288269043
$ws.Range("D1").Value = 1.339415431022644
$ws.Range("E1").Value = 0.7466664910316467
